$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) is treated as text so values like "1.000" or
# "0.00001204" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "23.512.60"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.649.22"
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "300.02"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("D7").Value = "0.3798"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.3575"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "0.08109"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "1.227"
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "6.426"
$ws.Range("E14").Value = "  -1.56%  "
$ws.Range("D15").Value = "7.430"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "0.00001204"
$ws.Range("E16").Value = "  -1.96%  "
$ws.Range("D17").Value = "1.654.85"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "97.33"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "0.06995"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "6.816"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "17.49"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").Value = "12.61"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "23.537.10"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").Value = "2.474"
$ws.Range("E25").Value = "  -1.93%  "
$ws.Range("D26").Value = "2.916"
$ws.Range("E26").Value = "  -6.45%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").Value = "153.30"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").Value = "5.235"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").Value = "133.35"
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("D31").Value = "1.840.98"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").Value = "6.948"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").Value = "2.142"
$ws.Range("E33").Value = "  +4.80%  "
$ws.Range("D34").Value = "11.93"
$ws.Range("E34").Value = "  +3.68%  "
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("D36").Value = "0.02737"
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("D37").Value = "0.08743"
$ws.Range("E37").Value = "  -0.72%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2457"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "6.001"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "13.31"
$ws.Range("E40").Value = "  +4.07%  "
$ws.Range("D41").Value = "0.06877"
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").Value = "0.6932"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "1.324"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").Value = "15.73"
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("E45").Value = "  -0.77%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "2.276"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "0.07819"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "128.13"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("D51").Value = "1.173"
$ws.Range("E51").Value = "  -1.23%  "
